$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New admission records to append (rows 96-108).
# Columns: A=F.I.SH, B=Ta'lim yo'nalishi, C=Passport, D=Shartnoma raqam,
#          E=Viloyat, F=Tuman, G=Telefon raqam, H=Sana
$data = @(
    @("Abdullaeva Xilola Ilhomovna", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AB0449918", "692", "Xorazm viloyati", "Xiva tumani", "998990610110", "01-11-2024"),
    @("Qurbonova Mohigul Esonali qizi", "Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik", "AD4463924", "693", "Fargona viloyati", "Beshariq tumani", "998948263202", "02-11-2024"),
    @("Djabborova Rushana Odil qizi", "Defektologiya (logopediya) 576 soatlik", "AB5761267", "694", "Buxoro viloyati", "Gʻijduvon tumani", "998916490504", "02-11-2024"),
    @("Sharibova Madina Farxodjon qizi", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AD4491495", "695", "Andijon viloyati", "Qoʻrgʻontepa tumani", "998996909708", "02-11-2024"),
    @("Nazirbayeva Dildora Xamza qizi", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AD3665614", "696", "Xorazm viloyati", "Gurlan tumani", "998972021510", "02-11-2024"),
    @("Kenjayeva Firuzabonu Avazbek qizi", "Defektologiya (logopediya) 576 soatlik", "AC1274290", "697", "Andijon viloyati", "Izboskan tuman", "998337200907", "02-11-2024"),
    @("Qayumova Xolidaxon Shuxratovna", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AD5502952", "698", "Fargona viloyati", "Buvayda tumani", "998900555522", "02-11-2024"),
    @("Kamolova Rushana Jamilovna", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AD4512393", "699", "Toshkent viloyati", "Ohangaron tumani", "998909481022", "03-11-2024"),
    @("Xudayberdiyeva Fotima Oybek qizi", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AB6207982", "700", "Qashqadaryo viloyati", "Nishon tumani", "998908747305", "03-11-2024"),
    @("Hamidova Muxtabar Obidovna", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AB2567833", "701", "Navoiy viloyati", "Tomdi tumani", "998936631068", "03-11-2024"),
    @("Kurbonova Feruza Baxtiyor qizi", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AA9126315", "702", "Samarqand viloyati", "Kattaqoʻrgʻon tumani", "998939958881", "04-11-2024"),
    @("Ismoilova Nozima Alisherovna", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AD2465041", "703", "Toshkent shahri", "Yashnaobod tumani", "998909343336", "04-11-2024"),
    @("Umarova Qumriyahon Qurbonovna", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AD8990709", "704", "Navoiy viloyati", "Konimex tumani", "998934311981", "04-11-2024")
)

$startRow = 96

# Columns whose values look numeric/date-like and would otherwise be
# auto-converted by Excel into real numbers/dates. We force them to stay
# plain text (matching the source data, which stores everything as text)
# by entering them with a leading apostrophe (text-prefix), then clearing
# the resulting cell format so no stray number formatting is left behind.
$textForcedCols = @(4, 7, 8)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    for ($col = 1; $col -le 8; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $value = $rec[$col - 1]
        if ($textForcedCols -contains $col) {
            $cell.Value = "'" + $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
